# Fruta / hortaliza, semanal
# Update the weekly price table: dates in column D and the associated
# volume / price figures in columns M, N, O, P, S are refreshed for each
# data row (rows 2-18) of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> @(Fecha(D), Volumen(M), PrecioMinimo(N), PrecioMaximo(O), PrecioPromedio(P), PrecioKg(S))
$rowData = @{
    2  = @(44895, 240, 3000,  3500,  3250,  1625)
    3  = @(44517, 400, 5500,  6000,  5750,  2875)
    4  = @(44818, 200, 11000, 12000, 11500, 5750)
    5  = @(44455, 200, 12000, 13000, 12500, 6250)
    6  = @(44475, 240, 11000, 12000, 11500, 5750)
    7  = @(44490, 400, 9500,  10000, 9750,  4875)
    8  = @(44875, 400, 7000,  7500,  7250,  3625)
    9  = @(44889, 460, 3500,  4000,  3750,  1875)
    10 = @(44881, 440, 6000,  7000,  6500,  3250)
    11 = @(44874, 300, 7500,  8000,  7750,  3875)
    12 = @(44461, 200, 11000, 12000, 11500, 5750)
    13 = @(44882, 440, 6000,  7000,  6500,  3250)
    14 = @(44482, 240, 10000, 11000, 10500, 5250)
    15 = @(44497, 500, 9000,  10000, 9500,  4750)
    16 = @(44454, 160, 12000, 13000, 12500, 6250)
    17 = @(44489, 160, 9500,  10000, 9750,  4875)
    18 = @(44819, 240, 11000, 12000, 11500, 5750)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]

    $ws.Cells.Item($row, 4).Value2  = $vals[0]   # D: Fecha
    $ws.Cells.Item($row, 13).Value2 = $vals[1]   # M: Volumen
    $ws.Cells.Item($row, 14).Value2 = $vals[2]   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value2 = $vals[3]   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value2 = $vals[4]   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value2 = $vals[5]   # S: Precio $/Kg
}
